# Atualização automática de RIO_PARDO.xlsx
#
# - Rename "Paineis DARQ"            -> "PAINEIS DARQ"
# - Rename "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO"
# - Delete "Desarquivamentos Pendentes" sheet entirely

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# Remove the obsolete sheet first.
$wb.Worksheets.Item("Desarquivamentos Pendentes").Delete()

# Normalize sheet names to upper case (with the correct accented
# characters) as requested by the update.
$wb.Worksheets.Item("Paineis DARQ").Name = "PAINEIS DARQ"
$wb.Worksheets.Item("Recolhimento x Eliminacao").Name = "RECOLHIMENTO X ELIMINAÇÃO"

# Deleting the trailing sheet moves Excel's active-tab focus; restore it
# back to the original selected sheet ("Paineis DARQ" / now "PAINEIS DARQ").
$wb.Worksheets.Item("PAINEIS DARQ").Activate()
